$wb = $excel.ActiveWorkbook

# --- Sheet "First run" ---
$ws1 = $wb.Worksheets.Item("First run")

# New header columns for TLD / Country
$ws1.Cells.Item(1, 8).Value = "TLD"
$ws1.Cells.Item(1, 9).Value = "Country"

# Row 5 (Google): Resource Type "N, o, n, e" -> "None"
$ws1.Range("F5").Value = "None"

# --- Sheet "Second run" ---
$ws2 = $wb.Worksheets.Item("Second run")

# New header columns for TLD / Country
$ws2.Cells.Item(1, 8).Value = "TLD"
$ws2.Cells.Item(1, 9).Value = "Country"

# Row 2 (Pangea / PANGAEA): add Country = European Union
$ws2.Range("I2").Value = "European Union"

# Rows 2-3: Resource Type "Data Center" -> "Catalog"
$ws2.Range("F2:F3").Value = "Catalog"

# Rows 11-25: Resource Type "Community" -> "Community, Image Collection, Organization"
$ws2.Range("F11:F25").Value = "Community, Image Collection, Organization"

# Remove the auto-generated hyperlinks (C3, C8) and restore plain formatting
$ws2.Hyperlinks.Delete()
$ws2.Range("C3").Style = "Normal"
$ws2.Range("C8").Style = "Normal"
